$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 7).Value  = -22.70742358078603   # G
    $ws.Cells.Item($r, 8).Value  = -22.70742358078603   # H
    $ws.Cells.Item($r, 9).Value  = -23.66812227074236   # I
    $ws.Cells.Item($r, 10).Value = -23.66812227074236   # J
    $ws.Cells.Item($r, 11).Value = -6.83                # K
    $ws.Cells.Item($r, 12).Value = -29.82532751091703   # L

    $ws.Cells.Item($r, 21).Value = 9.65                 # U
    $ws.Cells.Item($r, 22).Value = 1.655231560891938    # V
    $ws.Cells.Item($r, 23).Value = -0.04665300546448087 # W
    $ws.Cells.Item($r, 24).Value = 0.3908048681382948   # X
    $ws.Cells.Item($r, 25).Value = -0.4374578736027757  # Y
    $ws.Cells.Item($r, 26).Value = 0.001267083494715874 # Z
    $ws.Cells.Item($r, 27).Value = -0.02998948708017484 # AA
    $ws.Cells.Item($r, 28).Value = 0.06405757947912295  # AB
    $ws.Cells.Item($r, 29).Value = -0.0940470665592978  # AC
    $ws.Cells.Item($r, 30).Value = 55.9                 # AD
    $ws.Cells.Item($r, 31).Value = 0                    # AE
    $ws.Cells.Item($r, 32).Value = 55.9                 # AF
    $ws.Cells.Item($r, 33).Value = 46.25                # AG
    $ws.Cells.Item($r, 34).Value = 0.9055564555321561   # AH
    $ws.Cells.Item($r, 35).Value = 0.3029810298102981   # AI
    $ws.Cells.Item($r, 36).Value = 0.8880568356374808   # AJ
    $ws.Cells.Item($r, 37).Value = 0.2645124392336288   # AK
    $ws.Cells.Item($r, 38).Value = 1.49                 # AL
    $ws.Cells.Item($r, 39).Value = 1.424                # AM
    $ws.Cells.Item($r, 40).Value = -13.63414634146342   # AN
    $ws.Cells.Item($r, 41).Value = -3.63758389261745    # AO
    $ws.Cells.Item($r, 42).Value = -11.28048780487805   # AP
    $ws.Cells.Item($r, 43).Value = -3.806179775280899   # AQ
}
